$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Delete the "Instrument Detail" block (old rows 11-17): ---
#     Instrument Detail:, Spectrometer serial number, Grating Number,
#     Collection lens, Longpass filter, Entrance slit aperture, Fiber optic cable type
#     This shifts the remaining rows (Number of measurements.. Comments) up to rows 11-15,
#     which already land on the correct text/value/style without further edits.
$ws.Rows("11:17").Delete()

# --- 2. Row 1: clear the "Comment" value out of B1 (label stays "Parameter") ---
$ws.Range("B1").ClearContents()

# Update the help text in D1 to use the 12pt Calibri font (was 11pt)
$ws.Range("D1").Value = "The following information can be entered in the form. Information with a green background is required, others are optional."
$ws.Range("D1").Font.Size = 12

# --- 3. Row 2: Trial Name value + new helper note ---
$ws.Range("B2").Value = "CAP-2_2012_Aberdeen"
$ws.Range("D2").Value = 'This entry is defined in "Trial Name" of the "Phenotype Experiment"'

# --- 4. Row 6: relabel "start time" + add helper note ---
$ws.Range("A6").Value = "whole day, start time (24-hour clock)"
$ws.Range("D6").Value = "This will be ignored if the time is included in the raw data file"

# --- 5. Row 7: relabel "end time" + add helper note ---
$ws.Range("A7").Value = "whole day, end time (24-hour clock)"
$ws.Range("D7").Value = "This will be ignored if the time is included in the raw data file"

# --- 6. Row 8: add helper note for Integration time ---
$ws.Range("D8").Value = "This will be ignored if the integration time is included in the raw data file"

# --- 7. Row 10: Instrument -> Spectrometer System ---
$ws.Range("A10").Value = "Spectrometer System"
$ws.Range("B10").Value = "UCD_WUEoptimzed_Channel1"
$ws.Range("D10").Value = 'This entry is defined in the "System Name" of the "CSR Spectrometer System"'

# --- 8. Row 13 (was old row 20, "Reference (barium sulfate...)"): add unit note ---
$ws.Range("D13").Value = "s"

# --- 9. Selection moves from B5 to B2 ---
$ws.Range("B2").Select()

# --- 10. Column A width narrows from 65.33 to 50.5 characters ---
$ws.Columns(1).ColumnWidth = 49.666666666666664
